$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Cells.Item(53, 8).Value = 7225.3887
$ws.Cells.Item(53, 10).Value = 13297
$ws.Cells.Item(53, 12).Value = 13297
$ws.Cells.Item(53, 14).Value = -14571

# Row 64
$ws.Cells.Item(64, 8).Value = 5336.125
$ws.Cells.Item(64, 9).Value = 4789.3335
$ws.Cells.Item(64, 10).Value = 5664.2
$ws.Cells.Item(64, 11).Value = 4789.3335
$ws.Cells.Item(64, 12).Value = 5664.2
$ws.Cells.Item(64, 13).Value = -4541.3335
$ws.Cells.Item(64, 14).Value = -6160.2

# Row 67
$ws.Cells.Item(67, 8).Value = 5336.125
$ws.Cells.Item(67, 9).Value = 4789.3335
$ws.Cells.Item(67, 10).Value = 5664.2
$ws.Cells.Item(67, 11).Value = 4789.3335
$ws.Cells.Item(67, 12).Value = 5664.2
$ws.Cells.Item(67, 13).Value = -3931.3335
$ws.Cells.Item(67, 14).Value = -7380.2

# Row 74
$ws.Cells.Item(74, 8).Value = 3867.6
$ws.Cells.Item(74, 9).Value = 2334.5
$ws.Cells.Item(74, 10).Value = 10000
$ws.Cells.Item(74, 11).Value = 2334.5
$ws.Cells.Item(74, 12).Value = 10000
$ws.Cells.Item(74, 13).Value = -1398.5
$ws.Cells.Item(74, 14).Value = -11872

# Row 77
$ws.Cells.Item(77, 8).Value = 3867.6
$ws.Cells.Item(77, 9).Value = 2334.5
$ws.Cells.Item(77, 10).Value = 10000
$ws.Cells.Item(77, 11).Value = 11672.5
$ws.Cells.Item(77, 12).Value = 50000
$ws.Cells.Item(77, 13).Value = -6992.5
$ws.Cells.Item(77, 14).Value = -59360

# Row 138
$ws.Cells.Item(138, 8).Value = 3535.224
$ws.Cells.Item(138, 9).Value = 1187.3334
$ws.Cells.Item(138, 11).Value = 3562.0002
$ws.Cells.Item(138, 13).Value = 1577.9998


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 4466863
$ws.Cells.Item(32, 9).Value = 5001610.5
$ws.Cells.Item(32, 11).Value = 5001610.5
$ws.Cells.Item(32, 13).Value = -5001323.5

# Row 45
$ws.Cells.Item(45, 8).Value = 2677.3333
$ws.Cells.Item(45, 9).Value = 1725
$ws.Cells.Item(45, 11).Value = 1725
$ws.Cells.Item(45, 13).Value = -1348

# Row 97
$ws.Cells.Item(97, 8).Value = 4012
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 4012
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 4012
$ws.Cells.Item(97, 13).ClearContents()
$ws.Cells.Item(97, 14).Value = -5004

# Row 110
$ws.Cells.Item(110, 8).Value = 4510.2163
$ws.Cells.Item(110, 9).Value = 5583.1113
$ws.Cells.Item(110, 11).Value = 5583.1113
$ws.Cells.Item(110, 13).Value = -3538.1113

# Row 132
$ws.Cells.Item(132, 8).Value = 6967604.5
$ws.Cells.Item(132, 9).Value = 1650.5
$ws.Cells.Item(132, 11).Value = 4951.5
$ws.Cells.Item(132, 13).Value = -2421.5

# Row 133
$ws.Cells.Item(133, 8).Value = 69992
$ws.Cells.Item(133, 10).Value = 69992
$ws.Cells.Item(133, 12).Value = 69992
$ws.Cells.Item(133, 14).Value = -75052


$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 39272.523
$ws.Cells.Item(20, 9).Value = 17216.666
$ws.Cells.Item(20, 10).Value = 48094.867
$ws.Cells.Item(20, 11).Value = 17216.666
$ws.Cells.Item(20, 12).Value = 48094.867
$ws.Cells.Item(20, 13).Value = -16969.666
$ws.Cells.Item(20, 14).Value = -48588.867

# Row 86
$ws.Cells.Item(86, 8).Value = 1137.1
$ws.Cells.Item(86, 9).Value = 1174.5555
$ws.Cells.Item(86, 11).Value = 1174.5555
$ws.Cells.Item(86, 13).Value = -51.55549999999994

# Row 89
$ws.Cells.Item(89, 8).Value = 1137.1
$ws.Cells.Item(89, 9).Value = 1174.5555
$ws.Cells.Item(89, 11).Value = 5872.7775
$ws.Cells.Item(89, 13).Value = -256.7775000000001

# Row 94
$ws.Cells.Item(94, 8).Value = 1236.7916
$ws.Cells.Item(94, 9).Value = 1341.7894
$ws.Cells.Item(94, 11).Value = 1341.7894
$ws.Cells.Item(94, 13).Value = -890.7893999999999

# Row 97
$ws.Cells.Item(97, 8).Value = 926.5
$ws.Cells.Item(97, 9).Value = 926.5
$ws.Cells.Item(97, 11).Value = 926.5
$ws.Cells.Item(97, 13).Value = 64.5

# Row 107
$ws.Cells.Item(107, 8).Value = 1380.1666
$ws.Cells.Item(107, 9).Value = 1170.375
$ws.Cells.Item(107, 10).Value = 1799.75
$ws.Cells.Item(107, 11).Value = 1170.375
$ws.Cells.Item(107, 12).Value = 1799.75
$ws.Cells.Item(107, 13).Value = 749.625
$ws.Cells.Item(107, 14).Value = -5639.75


$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 5360.875
$ws.Cells.Item(16, 9).Value = 848
$ws.Cells.Item(16, 10).Value = 7412.1816
$ws.Cells.Item(16, 11).Value = 848
$ws.Cells.Item(16, 12).Value = 7412.1816
$ws.Cells.Item(16, 13).Value = -561
$ws.Cells.Item(16, 14).Value = -7986.1816

# Row 62
$ws.Cells.Item(62, 8).Value = 6537.25
$ws.Cells.Item(62, 9).Value = 6119.8
$ws.Cells.Item(62, 11).Value = 6119.8
$ws.Cells.Item(62, 13).Value = -5495.8

# Row 65
$ws.Cells.Item(65, 8).Value = 6537.25
$ws.Cells.Item(65, 9).Value = 6119.8
$ws.Cells.Item(65, 11).Value = 30599
$ws.Cells.Item(65, 13).Value = -27479

# Row 107
$ws.Cells.Item(107, 8).Value = 1394.8182
$ws.Cells.Item(107, 9).Value = 1183.875
$ws.Cells.Item(107, 10).Value = 1957.3334
$ws.Cells.Item(107, 11).Value = 1183.875
$ws.Cells.Item(107, 12).Value = 1957.3334
$ws.Cells.Item(107, 13).Value = 736.125
$ws.Cells.Item(107, 14).Value = -5797.3334

# Row 113
$ws.Cells.Item(113, 8).Value = 5360.875
$ws.Cells.Item(113, 9).Value = 848
$ws.Cells.Item(113, 10).Value = 7412.1816
$ws.Cells.Item(113, 11).Value = 848
$ws.Cells.Item(113, 12).Value = 7412.1816
$ws.Cells.Item(113, 13).Value = 1322
$ws.Cells.Item(113, 14).Value = -11752.1816


$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Cells.Item(52, 8).Value = 20717.875
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 10).Value = 20717.875
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 12).Value = 20717.875
$ws.Cells.Item(52, 13).ClearContents()
$ws.Cells.Item(52, 14).Value = -21235.875

# Row 70
$ws.Cells.Item(70, 8).Value = 8968.4
$ws.Cells.Item(70, 9).Value = 11643
$ws.Cells.Item(70, 10).Value = 4956.5
$ws.Cells.Item(70, 11).Value = 11643
$ws.Cells.Item(70, 12).Value = 4956.5
$ws.Cells.Item(70, 13).Value = -11373
$ws.Cells.Item(70, 14).Value = -5496.5

# Row 73
$ws.Cells.Item(73, 8).Value = 8968.4
$ws.Cells.Item(73, 9).Value = 11643
$ws.Cells.Item(73, 10).Value = 4956.5
$ws.Cells.Item(73, 11).Value = 11643
$ws.Cells.Item(73, 12).Value = 4956.5
$ws.Cells.Item(73, 13).Value = -10707
$ws.Cells.Item(73, 14).Value = -6828.5

# Row 121
$ws.Cells.Item(121, 8).Value = 25001
$ws.Cells.Item(121, 9).Value = 25001
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 11).Value = 25001
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 13).Value = -23254
$ws.Cells.Item(121, 14).ClearContents()


$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Cells.Item(46, 8).Value = 3780.5625
$ws.Cells.Item(46, 9).Value = 4998.5
$ws.Cells.Item(46, 10).Value = 3606.5715
$ws.Cells.Item(46, 11).Value = 4998.5
$ws.Cells.Item(46, 12).Value = 3606.5715
$ws.Cells.Item(46, 13).Value = -4810.5
$ws.Cells.Item(46, 14).Value = -3982.5715

# Row 74
$ws.Cells.Item(74, 8).Value = 55190
$ws.Cells.Item(74, 9).Value = 47975
$ws.Cells.Item(74, 10).Value = 60000
$ws.Cells.Item(74, 11).Value = 47975
$ws.Cells.Item(74, 12).Value = 60000
$ws.Cells.Item(74, 13).Value = -46977
$ws.Cells.Item(74, 14).Value = -61996

# Row 77
$ws.Cells.Item(77, 8).Value = 55190
$ws.Cells.Item(77, 9).Value = 47975
$ws.Cells.Item(77, 10).Value = 60000
$ws.Cells.Item(77, 11).Value = 143925
$ws.Cells.Item(77, 12).Value = 180000
$ws.Cells.Item(77, 13).Value = -138933
$ws.Cells.Item(77, 14).Value = -189984

# Row 82
$ws.Cells.Item(82, 8).Value = 1860.8823
$ws.Cells.Item(82, 9).Value = 1759.5
$ws.Cells.Item(82, 10).Value = 2334
$ws.Cells.Item(82, 11).Value = 1759.5
$ws.Cells.Item(82, 12).Value = 2334
$ws.Cells.Item(82, 13).Value = -1398.5
$ws.Cells.Item(82, 14).Value = -3056

# Row 85
$ws.Cells.Item(85, 8).Value = 1860.8823
$ws.Cells.Item(85, 9).Value = 1759.5
$ws.Cells.Item(85, 10).Value = 2334
$ws.Cells.Item(85, 11).Value = 1759.5
$ws.Cells.Item(85, 12).Value = 2334
$ws.Cells.Item(85, 13).Value = -511.5
$ws.Cells.Item(85, 14).Value = -4830

# Row 93
$ws.Cells.Item(93, 8).Value = 3511.65
$ws.Cells.Item(93, 9).Value = 6961.75
$ws.Cells.Item(93, 10).Value = 1211.5834
$ws.Cells.Item(93, 11).Value = 6961.75
$ws.Cells.Item(93, 12).Value = 1211.5834
$ws.Cells.Item(93, 13).Value = -5713.75
$ws.Cells.Item(93, 14).Value = -3707.5834

# Row 100
$ws.Cells.Item(100, 8).Value = 3049.9565
$ws.Cells.Item(100, 9).Value = 2634.7144
$ws.Cells.Item(100, 10).Value = 3695.889
$ws.Cells.Item(100, 11).Value = 2634.7144
$ws.Cells.Item(100, 12).Value = 3695.889
$ws.Cells.Item(100, 13).Value = -2093.7144
$ws.Cells.Item(100, 14).Value = -4777.889


$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 12799.042
$ws.Cells.Item(62, 10).Value = 10846.842
$ws.Cells.Item(62, 12).Value = 10846.842
$ws.Cells.Item(62, 14).Value = -12094.842

# Row 65
$ws.Cells.Item(65, 8).Value = 12799.042
$ws.Cells.Item(65, 10).Value = 10846.842
$ws.Cells.Item(65, 12).Value = 54234.21000000001
$ws.Cells.Item(65, 14).Value = -60474.21000000001

# Row 96
$ws.Cells.Item(96, 8).Value = 2000.75
$ws.Cells.Item(96, 10).Value = 2000.75
$ws.Cells.Item(96, 12).Value = 2000.75
$ws.Cells.Item(96, 14).Value = -4746.75

# Row 119
$ws.Cells.Item(119, 8).Value = 200000
$ws.Cells.Item(119, 10).Value = 200000
$ws.Cells.Item(119, 12).Value = 200000
$ws.Cells.Item(119, 14).Value = -209676

